# The workbook originally has a single sheet ("Sheet1") holding Thomas-model
# fitting results for several experiments. This edit:
#   1. Tidies a handful of now-unused "q_e" cells on Sheet1 (D7/D10/D11),
#      and normalises the number format on the D cells that still hold the
#      159.049 q_e value so they match the format used elsewhere in column D.
#   2. Updates the active selection on Sheet1.
#   3. Adds a new "Sheet2" (placed after Sheet1, becomes the active sheet)
#      containing the flattened/sorted flowrate + kTh (+ q_e) rows for every
#      experiment whose comment is "fixed q_e; min waste only" - this is the
#      data consumed by simulation_new_data.py per the commit message.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 touch-ups -------------------------------------------------

# D6/D9/D16/D21/D25 keep their value (159.049) but their number format is
# re-applied so it matches column D's other scientific-notation cells.
$ws1.Range("D6").NumberFormat  = "0.00E+00"
$ws1.Range("D9").NumberFormat  = "0.00E+00"
$ws1.Range("D16").NumberFormat = "0.00E+00"
$ws1.Range("D21").NumberFormat = "0.00E+00"
$ws1.Range("D25").NumberFormat = "0.00E+00"

# D7/D10/D11 were empty placeholder cells (no value, leftover formatting) -
# clear them out completely.
$ws1.Range("D7").Clear()
$ws1.Range("D10").Clear()
$ws1.Range("D11").Clear()

# Re-merge the expt_3/expt_4/expt_5 blocks so they sort after the
# expt_1/expt_2 merges in the saved merge-cell list.
$ws1.Range("A16:A19").UnMerge()
$ws1.Range("A16:A19").Merge()
$ws1.Range("A21:A23").UnMerge()
$ws1.Range("A21:A23").Merge()
$ws1.Range("E16:E19").UnMerge()
$ws1.Range("E16:E19").Merge()
$ws1.Range("E21:E23").UnMerge()
$ws1.Range("E21:E23").Merge()
$ws1.Range("A25:A27").UnMerge()
$ws1.Range("A25:A27").Merge()
$ws1.Range("E25:E27").UnMerge()
$ws1.Range("E25:E27").Merge()

# Selection on Sheet1 moves onto the expt_4 block.
$ws1.Range("B21:D27").Select()

# --- New Sheet2 ---------------------------------------------------------

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = 67.5
$ws2.Range("B1").Value = 0.776

$ws2.Range("A2").Value = 72.06
$ws2.Range("B2").Value = 0.4987
$ws2.Range("C2").Value = 159.049

$ws2.Range("A3").Value = 72.16
$ws2.Range("B3").Value = 0.216
$ws2.Range("C3").Value = 159.049

$ws2.Range("A4").Value = 72.16
$ws2.Range("B4").Value = 0.517

$ws2.Range("A5").Value = 74
$ws2.Range("B5").Value = 0.222006
$ws2.Range("C5").Value = 159.049

$ws2.Range("A6").Value = 74.67
$ws2.Range("B6").Value = 0.0716

$ws2.Range("A7").Value = 74.67
$ws2.Range("B7").Value = 0.00088

$ws2.Range("A8").Value = 75
$ws2.Range("B8").Value = 0.0024
$ws2.Range("C8").Value = 159.049

$ws2.Range("A9").Value = 75.5
$ws2.Range("B9").Value = 0.544281

$ws2.Range("A10").Value = 75.5
$ws2.Range("B10").Value = 0.6369
$ws2.Range("C10").Value = 159.049

$ws2.Range("A11").Value = 75.89
$ws2.Range("B11").Value = 0.000000576249698999999955397551

$ws2.Range("A12").Value = 77
$ws2.Range("B12").Value = 0.1969

$ws2.Range("A13").Value = 77.4
$ws2.Range("B13").Value = 0.4046

$ws2.Range("A14").Value = 77.5
$ws2.Range("B14").Value = 0.0021

$ws2.Range("A15").Value = 78
$ws2.Range("B15").Value = 0.00132799

# Columns B & C (kTh / q_e) use the same scientific-notation format as
# Sheet1's kTh column, except the very last row (B15) which stays General -
# matching the source data it was copied from (Sheet1!C11, unstyled).
$ws2.Range("B1:B14").NumberFormat = "0.00E+00"
$ws2.Range("C2").NumberFormat  = "0.00E+00"
$ws2.Range("C3").NumberFormat  = "0.00E+00"
$ws2.Range("C5").NumberFormat  = "0.00E+00"
$ws2.Range("C8").NumberFormat  = "0.00E+00"
$ws2.Range("C10").NumberFormat = "0.00E+00"

$ws2.Range("A1").Select()
